# Progress update as of 04-Nov-2025:
#  - Column H ("PERIOD TO EXPIRE") drops by 1 day for every data row.
#  - Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025.
# Both columns are plain "General"-formatted cells holding literal values
# (H = number, I = text that merely looks like a date), so we must avoid
# Excel's automatic "looks like a date" -> date-serial conversion when we
# write the new text into column I. We force a text entry with a leading
# apostrophe, then restore the cell's original (non quote-prefixed,
# General) formatting by pasting the number format from the sibling H
# cell on the same row, so the on-disk style stays exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow  = 40

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    # Column H: decrement the numeric "period to expire" by one day.
    $hCell.Value = $hCell.Value2 - 1

    # Column I: overwrite the literal date-text, forcing text (not a date
    # serial) via a leading apostrophe, then reapply the original General
    # number format (copied from column H on the same row) so the cell's
    # style stays unchanged.
    $iCell.Value = "'04-Nov-2025"
    $hCell.Copy()
    $iCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
